$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: swap GAS_NATURAL -> DIESEL, update date ---
$ws.Range("B2").Value = "DIESEL"
$ws.Range("E2").Value = 44856
$ws.Range("E2").NumberFormat = "mmm-yy"

# --- Rows 3-6: update date only (values otherwise unchanged) ---
$ws.Range("E3").Value = 44887
$ws.Range("E3").NumberFormat = "mmm-yy"
$ws.Range("E4").Value = 44887
$ws.Range("E4").NumberFormat = "mmm-yy"
$ws.Range("E5").Value = 44887
$ws.Range("E5").NumberFormat = "mmm-yy"
$ws.Range("E6").Value = 44887
$ws.Range("E6").NumberFormat = "mmm-yy"

# --- New rows 7-11: additional COMBUSTION_FIJA fuel entries ---
$ws.Range("A7").Value = "COMBUSTION_FIJA"
$ws.Range("B7").Value = "KEROSENE"
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = "MENSUAL"
$ws.Range("E7").Value = 44856
$ws.Range("E7").NumberFormat = "mmm-yy"

$ws.Range("A8").Value = "COMBUSTION_FIJA"
$ws.Range("B8").Value = "FUEL_OIL"
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = "MENSUAL"
$ws.Range("E8").Value = 44856
$ws.Range("E8").NumberFormat = "mmm-yy"

$ws.Range("A9").Value = "COMBUSTION_FIJA"
$ws.Range("B9").Value = "NAFTA"
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = "MENSUAL"
$ws.Range("E9").Value = 44856
$ws.Range("E9").NumberFormat = "mmm-yy"

$ws.Range("A10").Value = "COMBUSTION_FIJA"
$ws.Range("B10").Value = "CARBON_DE_LEÑA"
$ws.Range("C10").Value = 50
$ws.Range("D10").Value = "MENSUAL"
$ws.Range("E10").Value = 44856
$ws.Range("E10").NumberFormat = "mmm-yy"

$ws.Range("A11").Value = "COMBUSTION_FIJA"
$ws.Range("B11").Value = "LEÑA"
$ws.Range("C11").Value = 50
$ws.Range("D11").Value = "MENSUAL"
$ws.Range("E11").Value = 44856
$ws.Range("E11").NumberFormat = "mmm-yy"

# --- Column widths (auto-fit based widths observed in target) ---
$ws.Columns.Item(1).ColumnWidth = 31.666666666666668
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 8.666666666666666
$ws.Columns.Item(5).ColumnWidth = 8.5

# --- Selection cell matches target view state ---
$ws.Range("G12").Select()
